$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Rule R30's "From" threshold (C1) was corrected from 18 to 1.
$ws.Range("C10").Value = 1
